{"js": "// Replace arithmetic expressions (e.g. \"99-77=\") with new ones in the\n// worksheet table. The table has 20 rows x 5 columns = 100 cells, and the\n// replacement pairs below are listed in the same row-major order as the\n// cells in the document, so replacements[i] corresponds to cell i.\n// We resolve each target cell explicitly (rather than doing a single\n// whole-document text search) because some of the new expressions are\n// substrings of other new expressions (e.g. \"9+31=\" is contained in\n// \"39+31=\"), which would make a naive sequential body-wide search pick up\n// the wrong occurrence once earlier replacements have already happened.\nconst replacements = [\n  [\"99-77=\", \"65-33=\"],\n  [\"60-36=\", \"68+11=\"],\n  [\"54-50=\", \"59+34=\"],\n  [\"43-7=\", \"92-5=\"],\n  [\"65+6=\", \"53+6=\"],\n  [\"98-18=\", \"62-9=\"],\n  [\"20+31=\", \"64+33=\"],\n  [\"81-10=\", \"66-16=\"],\n  [\"91-8=\", \"29+47=\"],\n  [\"8+68=\", \"30+53=\"],\n  [\"48+17=\", \"81-65=\"],\n  [\"86-36=\", \"65-6=\"],\n  [\"61-61=\", \"92-85=\"],\n  [\"30-21=\", \"36+2=\"],\n  [\"73-37=\", \"16+47=\"],\n  [\"29+63=\", \"17+9=\"],\n  [\"39+16=\", \"99-3=\"],\n  [\"66-14=\", \"52-30=\"],\n  [\"4+94=\", \"87-41=\"],\n  [\"32+57=\", \"21+59=\"],\n  [\"17+10=\", \"33+56=\"],\n  [\"13+48=\", \"78-6=\"],\n  [\"69-62=\", \"36-27=\"],\n  [\"74-62=\", \"44+36=\"],\n  [\"16+77=\", \"43+52=\"],\n  [\"11+11=\", \"55+44=\"],\n  [\"64+23=\", \"56+11=\"],\n  [\"44-16=\", \"40+15=\"],\n  [\"75-28=\", \"57+11=\"],\n  [\"76+2=\", \"80-43=\"],\n  [\"81-23=\", \"0+21=\"],\n  [\"71-44=\", \"67+31=\"],\n  [\"33+59=\", \"35-22=\"],\n  [\"79-26=\", \"89-16=\"],\n  [\"16+67=\", \"66-56=\"],\n  [\"43+21=\", \"34+2=\"],\n  [\"18+31=\", \"6+19=\"],\n  [\"5+3=\", \"3+90=\"],\n  [\"13+2=\", \"48+2=\"],\n  [\"50-0=\", \"5+59=\"],\n  [\"66-63=\", \"68+9=\"],\n  [\"12+49=\", \"54-32=\"],\n  [\"48+50=\", \"36+34=\"],\n  [\"70-3=\", \"33+31=\"],\n  [\"94-53=\", \"3+90=\"],\n  [\"45+6=\", \"39+31=\"],\n  [\"90-4=\", \"38+25=\"],\n  [\"44-26=\", \"80-35=\"],\n  [\"93-59=\", \"15+65=\"],\n  [\"26-13=\", \"58-41=\"],\n  [\"84-54=\", \"22+50=\"],\n  [\"97-10=\", \"25+53=\"],\n  [\"21+17=\", \"64-18=\"],\n  [\"56-55=\", \"25-1=\"],\n  [\"62-17=\", \"29-21=\"],\n  [\"50-43=\", \"5+40=\"],\n  [\"46-28=\", \"73-64=\"],\n  [\"1+69=\", \"99-13=\"],\n  [\"76-28=\", \"48+14=\"],\n  [\"16+58=\", \"67-10=\"],\n  [\"84-37=\", \"99-81=\"],\n  [\"10+14=\", \"28+47=\"],\n  [\"91-69=\", \"32+21=\"],\n  [\"64-23=\", \"42-24=\"],\n  [\"16+80=\", \"45-34=\"],\n  [\"94-25=\", \"6+55=\"],\n  [\"31+21=\", \"95-35=\"],\n  [\"17-6=\", \"61-22=\"],\n  [\"0+68=\", \"90-9=\"],\n  [\"85-28=\", \"37+20=\"],\n  [\"8+8=\", \"79-38=\"],\n  [\"23+8=\", \"82-9=\"],\n  [\"81+5=\", \"90-18=\"],\n  [\"64-54=\", \"68-12=\"],\n  [\"75+24=\", \"54-33=\"],\n  [\"63-22=\", \"79+0=\"],\n  [\"21+50=\", \"17-7=\"],\n  [\"90-38=\", \"5+20=\"],\n  [\"35-3=\", \"57+11=\"],\n  [\"9+31=\", \"94-55=\"],\n  [\"86-50=\", \"39-28=\"],\n  [\"41+27=\", \"34+11=\"],\n  [\"16+70=\", \"54-47=\"],\n  [\"29-20=\", \"75-43=\"],\n  [\"83-35=\", \"83-39=\"],\n  [\"47-16=\", \"77-1=\"],\n  [\"26+8=\", \"53-0=\"],\n  [\"18+24=\", \"22+24=\"],\n  [\"36+1=\", \"50+17=\"],\n  [\"22-20=\", \"83-43=\"],\n  [\"27+12=\", \"60+26=\"],\n  [\"72-51=\", \"55+24=\"],\n  [\"23+44=\", \"89-47=\"],\n  [\"47+49=\", \"13-5=\"],\n  [\"6+71=\", \"1+82=\"],\n  [\"65-43=\", \"32-7=\"],\n  [\"75+13=\", \"78-74=\"],\n  [\"8+57=\", \"3+79=\"],\n  [\"43-34=\", \"56-11=\"],\n  [\"74+0=\", \"74+1=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Load the cells collection for every row up front.\nfor (const row of rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\nlet index = 0;\nfor (const row of rows.items) {\n  for (const cell of row.cells.items) {\n    if (index >= replacements.length) break;\n    const [oldText, newText] = replacements[index];\n\n    const results = cell.body.search(oldText, {matchCase: true, matchWholeWord: false});\n    results.load(\"items\");\n    await context.sync();\n\n    if (results.items.length > 0) {\n      results.items[0].insertText(newText, Word.InsertLocation.replace);\n    }\n\n    index++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace arithmetic expressions (e.g. \"99-77=\") with new ones in the\n# worksheet table. The table has 20 rows x 5 columns = 100 cells, laid out\n# row-major; $replacements lists the [old, new] pair for each cell in that\n# same row-major order, so $replacements[i] corresponds to cell i.\n# We address each cell explicitly by (row, column) via $t.Cell(r, c) and\n# assign Range.Text directly; this preserves each cell's paragraph/run\n# formatting (font, size, alignment) exactly like Word does for a plain\n# text replacement, and avoids any ambiguity from searching for text\n# fragments that might also appear inside newly-inserted replacement text.\n$replacements = @(\n    @(\"99-77=\", \"65-33=\"),\n    @(\"60-36=\", \"68+11=\"),\n    @(\"54-50=\", \"59+34=\"),\n    @(\"43-7=\", \"92-5=\"),\n    @(\"65+6=\", \"53+6=\"),\n    @(\"98-18=\", \"62-9=\"),\n    @(\"20+31=\", \"64+33=\"),\n    @(\"81-10=\", \"66-16=\"),\n    @(\"91-8=\", \"29+47=\"),\n    @(\"8+68=\", \"30+53=\"),\n    @(\"48+17=\", \"81-65=\"),\n    @(\"86-36=\", \"65-6=\"),\n    @(\"61-61=\", \"92-85=\"),\n    @(\"30-21=\", \"36+2=\"),\n    @(\"73-37=\", \"16+47=\"),\n    @(\"29+63=\", \"17+9=\"),\n    @(\"39+16=\", \"99-3=\"),\n    @(\"66-14=\", \"52-30=\"),\n    @(\"4+94=\", \"87-41=\"),\n    @(\"32+57=\", \"21+59=\"),\n    @(\"17+10=\", \"33+56=\"),\n    @(\"13+48=\", \"78-6=\"),\n    @(\"69-62=\", \"36-27=\"),\n    @(\"74-62=\", \"44+36=\"),\n    @(\"16+77=\", \"43+52=\"),\n    @(\"11+11=\", \"55+44=\"),\n    @(\"64+23=\", \"56+11=\"),\n    @(\"44-16=\", \"40+15=\"),\n    @(\"75-28=\", \"57+11=\"),\n    @(\"76+2=\", \"80-43=\"),\n    @(\"81-23=\", \"0+21=\"),\n    @(\"71-44=\", \"67+31=\"),\n    @(\"33+59=\", \"35-22=\"),\n    @(\"79-26=\", \"89-16=\"),\n    @(\"16+67=\", \"66-56=\"),\n    @(\"43+21=\", \"34+2=\"),\n    @(\"18+31=\", \"6+19=\"),\n    @(\"5+3=\", \"3+90=\"),\n    @(\"13+2=\", \"48+2=\"),\n    @(\"50-0=\", \"5+59=\"),\n    @(\"66-63=\", \"68+9=\"),\n    @(\"12+49=\", \"54-32=\"),\n    @(\"48+50=\", \"36+34=\"),\n    @(\"70-3=\", \"33+31=\"),\n    @(\"94-53=\", \"3+90=\"),\n    @(\"45+6=\", \"39+31=\"),\n    @(\"90-4=\", \"38+25=\"),\n    @(\"44-26=\", \"80-35=\"),\n    @(\"93-59=\", \"15+65=\"),\n    @(\"26-13=\", \"58-41=\"),\n    @(\"84-54=\", \"22+50=\"),\n    @(\"97-10=\", \"25+53=\"),\n    @(\"21+17=\", \"64-18=\"),\n    @(\"56-55=\", \"25-1=\"),\n    @(\"62-17=\", \"29-21=\"),\n    @(\"50-43=\", \"5+40=\"),\n    @(\"46-28=\", \"73-64=\"),\n    @(\"1+69=\", \"99-13=\"),\n    @(\"76-28=\", \"48+14=\"),\n    @(\"16+58=\", \"67-10=\"),\n    @(\"84-37=\", \"99-81=\"),\n    @(\"10+14=\", \"28+47=\"),\n    @(\"91-69=\", \"32+21=\"),\n    @(\"64-23=\", \"42-24=\"),\n    @(\"16+80=\", \"45-34=\"),\n    @(\"94-25=\", \"6+55=\"),\n    @(\"31+21=\", \"95-35=\"),\n    @(\"17-6=\", \"61-22=\"),\n    @(\"0+68=\", \"90-9=\"),\n    @(\"85-28=\", \"37+20=\"),\n    @(\"8+8=\", \"79-38=\"),\n    @(\"23+8=\", \"82-9=\"),\n    @(\"81+5=\", \"90-18=\"),\n    @(\"64-54=\", \"68-12=\"),\n    @(\"75+24=\", \"54-33=\"),\n    @(\"63-22=\", \"79+0=\"),\n    @(\"21+50=\", \"17-7=\"),\n    @(\"90-38=\", \"5+20=\"),\n    @(\"35-3=\", \"57+11=\"),\n    @(\"9+31=\", \"94-55=\"),\n    @(\"86-50=\", \"39-28=\"),\n    @(\"41+27=\", \"34+11=\"),\n    @(\"16+70=\", \"54-47=\"),\n    @(\"29-20=\", \"75-43=\"),\n    @(\"83-35=\", \"83-39=\"),\n    @(\"47-16=\", \"77-1=\"),\n    @(\"26+8=\", \"53-0=\"),\n    @(\"18+24=\", \"22+24=\"),\n    @(\"36+1=\", \"50+17=\"),\n    @(\"22-20=\", \"83-43=\"),\n    @(\"27+12=\", \"60+26=\"),\n    @(\"72-51=\", \"55+24=\"),\n    @(\"23+44=\", \"89-47=\"),\n    @(\"47+49=\", \"13-5=\"),\n    @(\"6+71=\", \"1+82=\"),\n    @(\"65-43=\", \"32-7=\"),\n    @(\"75+13=\", \"78-74=\"),\n    @(\"8+57=\", \"3+79=\"),\n    @(\"43-34=\", \"56-11=\"),\n    @(\"74+0=\", \"74+1=\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\n$index = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n  for ($c = 1; $c -le $colCount; $c++) {\n    if ($index -ge $replacements.Count) { break }\n\n    $pair = $replacements[$index]\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $cell = $t.Cell($r, $c)\n    $cellRange = $cell.Range\n    # Trim the trailing end-of-cell mark so Range.Text only covers the\n    # visible content, then overwrite it in place. Word keeps the\n    # surrounding run/paragraph formatting (font, size, alignment) when a\n    # range's text is reassigned this way.\n    $cellRange.MoveEnd(1, -1) | Out-Null\n    $cellRange.Text = $newText\n\n    $index++\n  }\n}\n"}
